# Daily scrape update - 2025-10-04 03:00:38 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data rows 2-13 with the freshly scraped opportunities ---
# Columns: A=ID, B=LINK, C=TITLE, D=COUNTRY, E=PREMIUM, F=APPLICANTS, G=DURATION, H=ORGANIZATION

$rows = @(
    @{ Id="1328264"; C="Customer Service"; D="Leiria, Portugal"; F="4 applicants"; G="9 - 12 Weeks"; H="OWHM Hospitality Management Lda" },
    @{ Id="1328227"; C="Policy & Advocacy Intern"; D="Hyderabad, Telangana, India"; F="0 applicants"; G="9 - 12 Weeks"; H="Arunodhaya Trust" },
    @{ Id="1328225"; C="Social Impact & Sustainability Intern"; D="Hyderabad, Telangana, India"; F="0 applicants"; G="9 - 12 Weeks"; H="Arunodhaya Trust" },
    @{ Id="1328098"; C="Arabic Translator"; D="Delhi, India"; F="0 applicants"; G="6 - 18 Months"; H="Pie Multilingual Services" },
    @{ Id="1327541"; C="Software Developer"; D="União das freguesias de Cascais e Estoril, Portugal"; F="84 applicants"; G="3 - 6 Months"; H="Dark Cloud" },
    @{ Id="1325377"; C="Business Development Intern"; D="Athens, Greece"; F="75 applicants"; G="9 - 12 Weeks"; H="Eutopians" },
    @{ Id="1323669"; C="interior designer"; D="Sousse, Tunisie"; F="19 applicants"; G="9 - 12 Weeks"; H="Zitouni Interior" },
    @{ Id="1323361"; C="Business Administration and Sustainable Marketing Intern"; D="Manipal, Karnataka, India"; F="17 applicants"; G="9 - 12 Weeks"; H="M.A.H.E." },
    @{ Id="1323077"; C="Pharmacy Intern"; D="Manipal, Karnataka, India"; F="11 applicants"; G="9 - 12 Weeks"; H="M.A.H.E." },
    @{ Id="1317306"; C="Web Developer"; D="8670 Aljezur, Portugal"; F="226 applicants"; G="9 - 12 Weeks"; H="Barranco Da Fonte" },
    @{ Id="1313548"; C="Electrical Engineering Intern"; D="Chandigarh, India"; F="14 applicants"; G="9 - 12 Weeks"; H="CGC technical Campus Jhanjeri" },
    @{ Id="1304736"; C="Guest Relations Manager"; D="Heraklion, Greece"; F="186 applicants"; G="9 - 12 Weeks"; H="Remarc Internation" }
)

$r = 2
foreach ($row in $rows) {
    # Opportunity ID must remain text (keep leading apostrophe so Excel stores it as a string)
    $ws.Cells.Item($r, 1).Value = "'" + $row.Id
    $ws.Cells.Item($r, 2).Value = "https://aiesec.org/opportunity/global-talent/" + $row.Id
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# --- Remove the three opportunities that fell out of today's scrape (old rows 14-16) ---
$ws.Range("A14:H16").EntireRow.Delete() | Out-Null

# --- Adjust column widths to fit the refreshed content ---
$ws.Columns.Item(3).ColumnWidth = 59 - 0.8333333333333333
$ws.Columns.Item(4).ColumnWidth = 54 - 0.8333333333333333
$ws.Columns.Item(6).ColumnWidth = 17 - 0.8333333333333333
$ws.Columns.Item(8).ColumnWidth = 34 - 0.8333333333333333
